# Fruta / hortaliza, semanal
# Insert a new weekly record at row 88 (pushing the existing rows 88-143
# down to 89-144) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 88:143 down one row, creating a blank row 88.
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new "Granada" record.
$ws.Range("A88").Value = 10
$ws.Range("B88").Value = "Vega Modelo de Temuco"
$ws.Range("C88").Value = "La Araucanía"
$ws.Range("D88").Value = 44767
$ws.Range("E88").Value = 9
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100104
$ws.Range("H88").Value = "Frutos de pepita"
$ws.Range("I88").Value = 100104001
$ws.Range("J88").Value = "Granada"
$ws.Range("K88").Value = "Wonderfull"
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 55
$ws.Range("N88").Value = 14000
$ws.Range("O88").Value = 14000
$ws.Range("P88").Value = 14000
$ws.Range("Q88").Value = "`$/bandeja 10 kilos granel"
$ws.Range("R88").Value = "Provincia de Limarí"
$ws.Range("S88").Value = 1400
$ws.Range("T88").Value = 10
